{"js": "// Add a blank \"Body Text\" paragraph and a dated \"Body Text\" paragraph\n// right after the opening image paragraph (before the \"Run Timing\n// Prospects for 2018\" heading).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document opens with a single paragraph that just holds the inline\n// masthead image (style \"FirstParagraph\"). Insert the new paragraphs\n// right after it.\nconst imageParagraph = paragraphs.items[0];\n\n// First new paragraph: empty, styled \"Body Text\".\nconst blankParagraph = imageParagraph.insertParagraph(\"\", \"After\");\nblankParagraph.style = \"Body Text\";\n\n// Second new paragraph: the draft date, also styled \"Body Text\".\nconst dateParagraph = blankParagraph.insertParagraph(\"Friday, June 1, 2018\", \"After\");\ndateParagraph.style = \"Body Text\";\n\nawait context.sync();\n", "ps1": "# Add a blank \"Body Text\" paragraph and a dated \"Body Text\" paragraph\n# right after the opening image paragraph (before the \"Run Timing\n# Prospects for 2018\" heading).\n\n$d = $word.ActiveDocument\n\n# The document opens with a single paragraph that just holds the inline\n# masthead image. Insert the new paragraphs right after it.\n$imageRange = $d.Paragraphs.Item(1).Range\n$imageRange.InsertParagraphAfter()\n\n# First new paragraph (index 2): empty, styled \"Body Text\".\n$blankRange = $d.Paragraphs.Item(2).Range\n$blankRange.Style = \"Body Text\"\n$blankRange.InsertParagraphAfter()\n\n# Second new paragraph (index 3): the draft date, also styled \"Body Text\".\n$dateRange = $d.Paragraphs.Item(3).Range\n$dateRange.Style = \"Body Text\"\n$dateRange.Text = \"Friday, June 1, 2018\"\n"}
